$wb = $excel.ActiveWorkbook
$wsProdutos = $wb.Worksheets.Item("produtos")
$wsMovimentos = $wb.Worksheets.Item("movimentos")

# Add new row 98 to "produtos" sheet
$wsProdutos.Cells.Item(98, 1).Value = 97
$wsProdutos.Cells.Item(98, 2).Value = "ee"
$wsProdutos.Cells.Item(98, 3).Value = ""
$wsProdutos.Cells.Item(98, 4).Value = "KG"
$wsProdutos.Cells.Item(98, 5).Value = 0

# Add new row 18 to "movimentos" sheet
$wsMovimentos.Cells.Item(18, 1).Value = 17
$wsMovimentos.Cells.Item(18, 2).Value = "Esponja Dupla Face"
$wsMovimentos.Cells.Item(18, 3).Value = "ENTRADA"
$wsMovimentos.Cells.Item(18, 4).Value = 2
$wsMovimentos.Cells.Item(18, 5).Value = "2026-01-06 15:06:26"
